$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAges = @(44.502796153293403, 43.249136385309498, 70.534044117337501, 60.297884429135898, 28.947498609165201, 44.237399945344201, 22.861358549721501, 26.8642004286618, 28.852769194930602, 37.397060158650703, 63.678600313806498, 38.274834169408102, 35.709057825661397, 35.252331416839802, 61.146731303258001, 41.902563974668404, 60.215272054777401, 48.670139556712897, 43.892403192284398, 30.68428214727, 38.733829335074297, 17.805307035613101, 53.034941997960601, 71.926019045827104, 33.005291999113702, 56.371913836741797, 61.4345337477247, 34.5678347668199, 26.782958586321499, 38.045568769629398, 34.309749164141799, 54.241569399891198, 78.062048429391297, 36.055116897033102, 51.964787928920998, 78.694055609234297, 30.922507672522698, 53.725952956756799, 41.172747004091001, 31.9800019324136, 60.893930916923402, 62.390872779190602, 59.874038541608897, 64.041462729193796, 44.624839745250902, 60.9951620238148, 47.8225942070098, 49.832205265994702, 42.0550146897126, 51.438039927463699, 43.995081399319503, 34.409583919236702, 54.888818347245603, 57.629854370446701, 44.258741955534902, 34.101684129671, 34.482167408219098, 59.696251658700596, 47.873132184310599, 43.715127086233402, 53.959985762769598, 55.359635911708899, 60.589206081275698, 60.9839408479017, 60.798889304574097, 72.7824200141971, 68.359729955379393, 62.429579555569902, 31.501989231636099, 68.050907676778095, 41.367795209537199, 72.304010441793395, 34.038401742966698, 62.215863974918797, 51.551479755282102, 38.562249532204802, 64.0528285786, 77.947837237945393, 39.556963664379701, 74.702294608375297, 70.587764463974494, 59.656884795798199, 57.0152469251261, 51.269466875791203, 71.984129879353603, 43.076026957197797, 37.556912414350201, 55.7041120253446, 20.943085863621899, 36.579051167564202, 40.840958274285597, 34.602965234672503, 53.8435450650815, 40.892901075632899, 42.357796406351298, 60.7971331625658, 43.481933475517202, 45.0168399689098, 62.480868392241703, 42.337182198276302)

for ($i = 0; $i -lt $newAges.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newAges[$i]
}

$ws.Range("F9").Select()
